$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    0.99999999446970245,
    0.9944998079261258,
    0.97214215219555355,
    0.96152676808084692,
    0.95129592019270204,
    0.92663871104989726,
    0.92236117510864024,
    0.91620857385910226,
    0.90981930224196872,
    0.90434502994904264,
    0.90355327550941966,
    0.90230268269265901,
    0.89876440425612647,
    0.8981998272798537,
    0.89886605953967713,
    0.89635948494246276,
    0.89265147740156348,
    0.89154255851653308,
    0.99209588929849657,
    0.97294572906710941,
    0.97154722220467449,
    0.9702827112085266,
    0.9651064573831275,
    0.95208466699319783,
    0.94562751775399945,
    0.92413376529779023,
    0.92089104295302016,
    0.9092254520101043,
    0.90125291546560948,
    0.89847003688903049,
    0.9000192021001816,
    0.90149801393925033,
    0.89937853159520675
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}
